# Applies the "Added last minute updates" edit to the first paragraph of
# the document:
#   - adds a 5pt paragraph border (top/left/bottom/right) with no line style
#   - changes the left indent from 120 (6pt) to 225 (11.25pt) twips
#   - renames the placeholder token from
#       **ID__AFFARS_5303_topic_12__ID**
#     to
#       **ID__AFFARS_SUBPART_5303_6__ID**
#     and drops the trailing space run that followed it.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# --- paragraph border: <w:pBdr><w:top w:space="5"/>... (no val/sz/color) ---
# Setting the Distance* properties directly (without first flipping
# Borders.Enable) yields a bare pBdr with only w:space set, matching the
# target markup.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- left indent: 120 twips (6pt) -> 225 twips (11.25pt) ---
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- placeholder text + drop the trailing " " run ---
# Replacing the full original text (token + trailing space) with the new
# token (no trailing space) both renames the placeholder and removes the
# separate space-only run in one edit.
$find = $p1.Range.Find
$find.Execute("**ID__AFFARS_5303_topic_12__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5303_6__ID**", 2) | Out-Null

Write-Host "Paragraph 1 now reads:" $p1.Range.Text
